$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New url/status data for rows 2-33 (row 1 is the header "url"/"status", unchanged)
$data = @(
    @("https://preview.abbviecare.fr", 200),
    @("https://preview.abbviecare.fr/bin/public/abbvie-commons/basic-login", 302),
    @("https://preview.abbviecare.fr/libs/granite/csrf/token.json", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-header-publish.min.css", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-main.min.css", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-components.min.css", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/components/content/button-link/v2/button-link/clientlibs.min.js", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-footer-publish.min.js", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-main.min.js", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-components.min.js", 200),
    @("https://consent.trustarc.com/v2/notice/0ivu3f", 200),
    @("https://preview.abbviecare.fr/", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/assets/resources/fonts/NeueHaasUnicaPro-Light.woff2", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/assets/resources/fonts/NeueHaasUnicaPro-Bold.woff2", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/assets/resources/fonts/NeueHaasUnicaPro-Regular.woff2", 200),
    @("https://consent.trustarc.com/v2/asset/ic-error.svg", 200),
    @("https://consent.trustarc.com/v2/asset/ic-close-white.svg", 200),
    @("https://consent.trustarc.com/v2/asset/ic-close.svg", 200),
    @("https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg", 200),
    @("https://consent.trustarc.com/v2/asset/16:19:48.8270ivu3f_AbbVieID-logo.png", 200),
    @("https://consent.trustarc.com/v2/asset/latin.woff2", 200),
    @("https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=0ivu3f&behavior=implied&country=bd&language=en&rand=0.29477825510117084&session=f73ab2cc-6b3b-4eb3-8280-502fd2fe0a4a&userType=NEW", 202),
    @("https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=0ivu3f&referer=&fullURL=https%3A%2F%2Fpreview.abbviecare.fr%2F&category=", 200),
    @("https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=0ivu3f&behavior=implied&country=bd&language=en&rand=0.7043335869565213&session=f73ab2cc-6b3b-4eb3-8280-502fd2fe0a4a&userType=NEW", 202),
    @("https://preview.abbviecare.fr/libs/granite/csrf/token.json", 200),
    @("https://preview.abbviecare.fr/content/dam/abbvie-care-ous/fr/design-assets/banner-desktop.jpg/_jcr_content/renditions/cq5dam.web.1280.1280.jpeg", 200),
    @("https://preview.abbviecare.fr/content/dam/abbvie-care-ous/fr/design-assets/banner_logo.png/_jcr_content/renditions/cq5dam.web.1280.1280.png", 200),
    @("https://preview.abbviecare.fr/etc.clientlibs/awcm-projects-ous/components/content/external-link-popup/v2/external-link-popup/clientlibs.min.js", 200),
    @("https://preview.abbviecare.fr/bin/public/abbvie-commons/hreflangs?resourcePath=/content/abbvie-care-ous/fr/fr/jcr:content", 200),
    @("https://preview.abbviecare.fr/content/abbvie-care-ous/fr/fr/jcr:content/body/banner/banner-content/column_control/par1-100col/image_extension/item_1.coreimg.png/1663712714254-banner_logo.png", 200),
    @("https://preview.abbviecare.fr/content/abbvie-care-ous/fr/fr/jcr:content/body/banner/bg-image/item_1.coreimg.jpg/1663712702019-banner-desktop.jpg", 200)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Rows 34 and 35 no longer exist in the new data (table shrank from 35 to 33 rows).
# Delete the now-unused trailing rows 34 and 35 so the sheet dimension becomes A1:B33.
$ws.Range("A34:B35").Delete()
